$d = $word.ActiveDocument
$range = $d.Content
$found = $range.Find.Execute("same row, column, or diagonal.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "Found: $found"
Write-Output "Range start/end: $($range.Start) $($range.End)"
$range.Collapse(0)
Write-Output "collapsed start/end: $($range.Start) $($range.End)"
$range.Select()
$sel = $word.Selection
Write-Output "sel after select(): $($sel.Start) $($sel.End)"
$sel.TypeText(" A brute force algorithm testing every possible grid alignment takes quite a while because the grid is 8 by 8, but the algorithm can be simplified by first preventing queens from taking the same row or column as the other. ")
